$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 90.666664
$ws.Range("I9").Value = 88.5
$ws.Range("K9").Value = 88.5
$ws.Range("M9").Value = 80.5
$ws.Range("H92").Value = 14706175
$ws.Range("I92").Value = 16666865
$ws.Range("K92").Value = 16666865
$ws.Range("M92").Value = -16665617
$ws.Range("H112").Value = 1619.2368
$ws.Range("J112").Value = 1656.4445
$ws.Range("L112").Value = 4969.333500000001
$ws.Range("N112").Value = -7185.333500000001
$ws.Range("H132").Value = 868.8431399999999
$ws.Range("I132").Value = 819.3261
$ws.Range("K132").Value = 2457.9783
$ws.Range("M132").Value = 72.02170000000024
$ws.Range("H137").Value = 112501.22
$ws.Range("I137").Value = 400
$ws.Range("K137").Value = 1200
$ws.Range("M137").Value = 1350
$ws.Range("H138").Value = 2086.8223
$ws.Range("I138").Value = 2220.3784
$ws.Range("J138").Value = 1993.585
$ws.Range("K138").Value = 6661.135200000001
$ws.Range("L138").Value = 5980.755
$ws.Range("M138").Value = -1521.135200000001
$ws.Range("N138").Value = -16260.755
$ws.Range("H139").Value = 59679.5
$ws.Range("J139").Value = 59679.5
$ws.Range("L139").Value = 59679.5
$ws.Range("N139").Value = -69959.5
$ws.Range("H141").Value = 2801756.8
$ws.Range("I141").Value = 3500876.5
$ws.Range("K141").Value = 10502629.5
$ws.Range("M141").Value = -10497449.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1260.5652
$ws.Range("I45").Value = 1003.36365
$ws.Range("J45").Value = 1496.3334
$ws.Range("K45").Value = 1003.36365
$ws.Range("L45").Value = 1496.3334
$ws.Range("M45").Value = -626.36365
$ws.Range("N45").Value = -2250.3334
$ws.Range("H61").Value = 10216.35
$ws.Range("I61").Value = 9152.5
$ws.Range("J61").Value = 11812.125
$ws.Range("K61").Value = 9152.5
$ws.Range("L61").Value = 11812.125
$ws.Range("M61").Value = -8940.5
$ws.Range("N61").Value = -12236.125
$ws.Range("H122").Value = 1345.4615
$ws.Range("I122").Value = 1249.25
$ws.Range("K122").Value = 3747.75
$ws.Range("M122").Value = -1297.75
$ws.Range("H136").Value = 10216.35
$ws.Range("I136").Value = 9152.5
$ws.Range("J136").Value = 11812.125
$ws.Range("K136").Value = 27457.5
$ws.Range("L136").Value = 35436.375
$ws.Range("M136").Value = -24907.5
$ws.Range("N136").Value = -40536.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 269.4
$ws.Range("I64").Value = 314.5
$ws.Range("K64").Value = 314.5
$ws.Range("M64").Value = -89.5
$ws.Range("H67").Value = 269.4
$ws.Range("I67").Value = 314.5
$ws.Range("K67").Value = 314.5
$ws.Range("M67").Value = 465.5
$ws.Range("H99").Value = 1382.6
$ws.Range("J99").Value = 1519
$ws.Range("L99").Value = 1519
$ws.Range("N99").Value = -4515
$ws.Range("H134").Value = 7772.952
$ws.Range("I134").Value = 9356.4375
$ws.Range("K134").Value = 28069.3125
$ws.Range("M134").Value = -25534.3125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2214.3333
$ws.Range("I31").Value = 1205
$ws.Range("J31").Value = 3021.8
$ws.Range("K31").Value = 1205
$ws.Range("L31").Value = 3021.8
$ws.Range("M31").Value = -910
$ws.Range("N31").Value = -3611.8
$ws.Range("H34").Value = 2214.3333
$ws.Range("I34").Value = 1205
$ws.Range("J34").Value = 3021.8
$ws.Range("K34").Value = 1205
$ws.Range("L34").Value = 3021.8
$ws.Range("M34").Value = -1003
$ws.Range("N34").Value = -3425.8
$ws.Range("H105").Value = 866.25
$ws.Range("I105").Value = 850.6667
$ws.Range("K105").Value = 850.6667
$ws.Range("M105").Value = 896.3333
$ws.Range("H122").Value = 3674.2222
$ws.Range("I122").Value = 2185.5
$ws.Range("J122").Value = 4865.2
$ws.Range("K122").Value = 6556.5
$ws.Range("L122").Value = 14595.6
$ws.Range("M122").Value = -4106.5
$ws.Range("N122").Value = -19495.6
$ws.Range("H134").Value = 2484.027
$ws.Range("I134").Value = 2315
$ws.Range("J134").Value = 4399.6665
$ws.Range("K134").Value = 6945
$ws.Range("L134").Value = 13198.9995
$ws.Range("M134").Value = -4410
$ws.Range("N134").Value = -18268.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 142957920
$ws.Range("I50").Value = 174966
$ws.Range("J50").Value = 333335200
$ws.Range("K50").Value = 524898
$ws.Range("L50").Value = 1000005600
$ws.Range("M50").Value = -524417
$ws.Range("N50").Value = -1000006562
$ws.Range("H53").Value = 142957920
$ws.Range("I53").Value = 174966
$ws.Range("J53").Value = 333335200
$ws.Range("K53").Value = 524898
$ws.Range("L53").Value = 1000005600
$ws.Range("M53").Value = -524417
$ws.Range("N53").Value = -1000006562
$ws.Range("H98").Value = 1005.4
$ws.Range("J98").Value = 1061.5555
$ws.Range("L98").Value = 3184.6665
$ws.Range("N98").Value = -6180.666499999999
$ws.Range("H122").Value = 692.65
$ws.Range("J122").Value = 807.9286
$ws.Range("L122").Value = 7271.3574
$ws.Range("N122").Value = -12171.3574
$ws.Range("H141").Value = 3978.6667
$ws.Range("I141").Value = 3297.5557
$ws.Range("K141").Value = 9892.667099999999
$ws.Range("M141").Value = -4712.667099999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 14000
$ws.Range("I26").Value = 14000
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 14000
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -13720
$ws.Range("N26").ClearContents()
$ws.Range("H50").Value = 14000
$ws.Range("I50").Value = 14000
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 14000
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -13502
$ws.Range("N50").ClearContents()
$ws.Range("H102").Value = 2218.0476
$ws.Range("I102").Value = 2221.0557
$ws.Range("K102").Value = 2221.0557
$ws.Range("M102").Value = -599.0556999999999
$ws.Range("H122").Value = 1533.8334
$ws.Range("I122").Value = 1067.6666
$ws.Range("K122").Value = 3202.9998
$ws.Range("M122").Value = -752.9998000000001
$ws.Range("H126").Value = 2980006.2
$ws.Range("I126").Value = 3971865.8
$ws.Range("J126").Value = 202799.6
$ws.Range("K126").Value = 11915597.4
$ws.Range("L126").Value = 608398.8
$ws.Range("M126").Value = -11913127.4
$ws.Range("N126").Value = -613338.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2587.6
$ws.Range("I7").Value = 2616.1667
$ws.Range("K7").Value = 2616.1667
$ws.Range("M7").Value = -2504.1667
$ws.Range("I61").Value = 4299.2856
$ws.Range("J61").Value = 4712.7144
$ws.Range("K61").Value = 4299.2856
$ws.Range("L61").Value = 4712.7144
$ws.Range("M61").Value = -4097.2856
$ws.Range("N61").Value = -5116.7144
$ws.Range("H100").Value = 1486.875
$ws.Range("I100").Value = 1442.1428
$ws.Range("J100").Value = 1800
$ws.Range("K100").Value = 1442.1428
$ws.Range("L100").Value = 1800
$ws.Range("M100").Value = -901.1428000000001
$ws.Range("N100").Value = -2882
$ws.Range("I113").Value = 4299.2856
$ws.Range("J113").Value = 4712.7144
$ws.Range("K113").Value = 4299.2856
$ws.Range("L113").Value = 4712.7144
$ws.Range("M113").Value = -2129.2856
$ws.Range("N113").Value = -9052.714400000001
$ws.Range("H126").Value = 2587.6
$ws.Range("I126").Value = 2616.1667
$ws.Range("K126").Value = 7848.500100000001
$ws.Range("M126").Value = -5378.500100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 7084.1
$ws.Range("I126").Value = 14626.375
$ws.Range("K126").Value = 43879.125
$ws.Range("M126").Value = -41409.125
